$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bad Drivers summary row (row 3) ---
$ws.Range("C3").Value = 169
$ws.Range("D3").Value = 94.3

# --- Bad Drivers totals row (row 4) ---
$ws.Range("C4").Value = 169

# --- Good Drivers table re-sorted / refreshed (rows 12-18) ---
# Row 12
$ws.Range("A12").Value = "Intel(R) Wi-Fi 6E AX210 160MHz - 22.0.1.5"
$ws.Range("B12").Value = 156943
$ws.Range("D12").Value = 100

# Row 13
$ws.Range("A13").Value = "Intel(R) Wi-Fi 6E AX210 160MHz - 23.120.0.3"
$ws.Range("B13").Value = 34181
$ws.Range("D13").Value = 99.90000000000001

# Row 14
$ws.Range("A14").Value = "Intel(R) Wi-Fi 6E AX210 160MHz - 23.20.1.1"
$ws.Range("B14").Value = 13533
$ws.Range("D14").Value = 100

# Row 15
$ws.Range("A15").Value = "Intel(R) Wi-Fi 6E AX210 160MHz - 22.170.2.1"
$ws.Range("B15").Value = 19083
$ws.Range("D15").Value = 100

# Row 16
$ws.Range("A16").Value = "Intel(R) Wi-Fi 6E AX210 160MHz - 22.100.0.3"
$ws.Range("B16").Value = 12988
$ws.Range("D16").Value = 100

# Row 17
$ws.Range("A17").Value = "Intel(R) Wi-Fi 6E AX210 160MHz - 22.130.0.5"
$ws.Range("B17").Value = 18738
$ws.Range("D17").Value = 99.90000000000001

# Row 18 (A18/B18/D18 unchanged - only vintage date changes below)

# --- Driver Vintage (column E) dates: write as literal text, not dates ---
# Format the cells as Text first so the date-looking strings are not
# auto-converted into date serial numbers, then restore the original
# "General" right-aligned look (style 4) by copying formats from an
# untouched style-4 cell (B3), which leaves the text value intact.
$ws.Range("E13:E18").NumberFormat = "@"

$ws.Range("E13").Value = "2025-02-05"
$ws.Range("E14").Value = "2023-12-19"
$ws.Range("E15").Value = "2022-08-30"
$ws.Range("E16").Value = "2022-05-01"
$ws.Range("E17").Value = "2022-03-14"
$ws.Range("E18").Value = "2022-01-01"

$ws.Range("B3").Copy()
$ws.Range("E13:E18").PasteSpecial(-4122)
